# Updated cryptos list data: applies the per-cell text updates described by the diff.
# Values are written via a NumberFormat('@') round-trip so that numeric-looking strings
# (e.g. "212.31", "0.500") are stored as plain text instead of being coerced into
# numbers, then ClearFormats() drops the temporary format so the cell keeps the
# workbook's original (unstyled) look - matching the source file's plain inline-string
# cells, which carry no explicit style either.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextValue 'D2' '26.390.39'
Set-TextValue 'D3' '1.606.91'
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '212.31'
Set-TextValue 'E5' '  -0.04%  '
Set-TextValue 'D6' '0.500'
Set-TextValue 'E6' '  -0.39%  '
Set-TextValue 'E8' '  -0.04%  '
Set-TextValue 'E9' '  -0.06%  '
Set-TextValue 'D10' '19.28'
Set-TextValue 'E10' '  +1.54%  '
Set-TextValue 'D11' '0.0856'
Set-TextValue 'E11' '  +0.42%  '
Set-TextValue 'D12' '1.833.65'
Set-TextValue 'E12' '  +1.01%  '
Set-TextValue 'D13' '1.601.98'
Set-TextValue 'E13' '  +0.98%  '
Set-TextValue 'E14' '  +0.04%  '
Set-TextValue 'D15' '0.508'
Set-TextValue 'E15' '  -0.23%  '
Set-TextValue 'D16' '63.43'
Set-TextValue 'E16' '  -0.61%  '
Set-TextValue 'D17' '26.389.61'
Set-TextValue 'E17' '  +0.58%  '
Set-TextValue 'D18' '232.23'
Set-TextValue 'E18' '  +7.73%  '
Set-TextValue 'D19' '7.69'
Set-TextValue 'E19' '  +5.25%  '
Set-TextValue 'E20' '  -0.22%  '
Set-TextValue 'E22' '  -0.33%  '
Set-TextValue 'E23' '  -0.78%  '
Set-TextValue 'E24' '  +1.85%  '
Set-TextValue 'D25' '147.19'
Set-TextValue 'E25' '  +1.77%  '
Set-TextValue 'E26' '  -0.02%  '
Set-TextValue 'D27' '6.97'
Set-TextValue 'E27' '  +0.05%  '
Set-TextValue 'E28' '  +1.21%  '
Set-TextValue 'D29' '15.46'
Set-TextValue 'E29' '  +2.25%  '
Set-TextValue 'E30' '  +1.15%  '
Set-TextValue 'E31' '  +0.00%  '
Set-TextValue 'B32' 'Maker'
Set-TextValue 'C32' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D32' '1.489.06'
Set-TextValue 'E32' '  +5.72%  '
Set-TextValue 'B33' 'Filecoin'
Set-TextValue 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '3.23'
Set-TextValue 'E33' '  +1.39%  '
Set-TextValue 'E34' '  -0.47%  '
Set-TextValue 'E35' '  -0.48%  '
Set-TextValue 'E36' '  +0.95%  '
Set-TextValue 'D37' '0.562'
Set-TextValue 'E37' '  -2.88%  '
Set-TextValue 'E38' '  -0.01%  '
Set-TextValue 'E39' '  +0.04%  '
Set-TextValue 'E40' '  +0.17%  '
Set-TextValue 'E41' '  +0.04%  '
Set-TextValue 'D43' '0.936'
Set-TextValue 'E43' '  -4.51%  '
Set-TextValue 'D44' '1.744.83'
Set-TextValue 'E44' '  +0.98%  '
Set-TextValue 'E45' '  -0.78%  '
Set-TextValue 'E46' '  +0.12%  '
Set-TextValue 'D47' '89.24'
Set-TextValue 'E47' '  +3.37%  '
Set-TextValue 'E48' '  +0.48%  '
Set-TextValue 'E49' '  +0.06%  '
Set-TextValue 'D50' '0.0962'
Set-TextValue 'E50' '  +0.84%  '
Set-TextValue 'B51' 'EnergySwap'
Set-TextValue 'C51' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '7.45'
Set-TextValue 'E51' '  +0.84%  '
